$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append a new run (own run, own rPr sz24/szCs24) to the paragraph
#    that ends with the "... hættur að vinna"" text, right before that
#    paragraph's end-of-paragraph mark.
# ---------------------------------------------------------------------
$p29 = $d.Paragraphs.Item(29)
$insPoint = $d.Range($p29.Range.End - 1, $p29.Range.End - 1)
$newText = ", Notast var við meira af Stored Procedures í gagnagrunninum, grunnurinn að vefsíðu hent upp í flýti og sett upp IIS og PHP þjónustu fyrir serverinn og viðeigandi DNS þannig að serverinn og vélar tengdar honum gera farið inná Keplergames.com og farið þaðan á heimasíðuna."
$startPos = $insPoint.Start
$insPoint.InsertAfter($newText)

# Give the freshly inserted text its own run formatting (sz=24/12pt).
$newRunRange = $d.Range($startPos, $startPos + $newText.Length)
$newRunRange.Font.Size = 12

# Font.Size alone only stamps <w:sz>; stamp <w:szCs> too (as in the
# surrounding runs) by re-asserting the run's OOXML for just that span.
$newRunXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRunRange2 = $d.Range($startPos, $startPos + $newText.Length)
$newRunRange2.InsertXML($newRunXml)

# ---------------------------------------------------------------------
# 2. The "_GoBack" bookmark currently lives alone in the next (empty)
#    paragraph. Merge that paragraph into the one above (delete the
#    paragraph mark that separates them) so the bookmark now sits at
#    the end of our paragraph, right after the text we just added.
# ---------------------------------------------------------------------
$p29b = $d.Paragraphs.Item(29)
$markRange = $d.Range($p29b.Range.End - 1, $p29b.Range.End)
$markRange.Delete()

# ---------------------------------------------------------------------
# 3. Re-create a fresh empty paragraph (same formatting: centered,
#    sz24/szCs24) where the bookmark paragraph used to be, now that the
#    bookmark itself has moved up into the previous paragraph.
# ---------------------------------------------------------------------
$p30 = $d.Paragraphs.Item(30)
$r30 = $p30.Range
$r30.Collapse(1)
$r30.InsertParagraphBefore()

# Clean up the stray empty run InsertParagraphBefore leaves behind so
# the new paragraph only carries paragraph-level formatting.
$pBlank = $d.Paragraphs.Item(30)
$blankRange = $pBlank.Range
$blankXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$blankRange.InsertXML($blankXml)

Write-Output "edit complete"
